# promo_test_data.xlsx edit
# - Update the test-data email used across all rows (B3:B69)
# - Update the "user_mobile"-style numeric id in column C (C3:C69): 123456 -> 1234
# - Reset the sheet's scroll position / selection (was topLeftCell A61 / D78 selected)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update email used for every test row ---
$ws.Range("B3:B69").Value = "Zwingautomation78@gmail.com"

# --- Update numeric id (user_mobile) for every test row ---
$ws.Range("C3:C69").Value = 1234

# --- Reset view: clear scrolled-away top-left cell and select C2 ---
$ws.Activate() | Out-Null
$ws.Range("C2").Select() | Out-Null
